# "Cuisines and Restaurants.xlsx" — updated the list so no blank space in cuisines
#
# The filtered Cuisines list in column C of the "Cuisines" sheet had a blank
# cell at C11 (between "Desserts" and "Fast Food"), leaving a gap in the
# alphabetical list used for the Restaurants sheet's data-validation dropdown.
# Close the gap by shifting every entry below the blank up by one row, which
# also pushes one extra blank cell onto the end of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cuisines")

# Values currently in C12:C30, in order - these are the ones that need to move
# up into C11:C29 once the blank at C11 is closed up.
$cuisines = @(
    "Fast Food",
    "Fruits & Vegetables",
    "Halal",
    "Healthy",
    "Indian",
    "Italian",
    "Japanese",
    "Korean",
    "Local",
    "Malay",
    "Mexican",
    "Peranakan",
    "Salad",
    "Seafood",
    "Snacks",
    "Soup",
    "Thai",
    "Vietnamese",
    "Western"
)

for ($i = 0; $i -lt $cuisines.Length; $i++) {
    $ws.Cells.Item(11 + $i, 3).Value = $cuisines[$i]
}

# The last populated row (previously C30) is now vacated by the shift.
$ws.Cells.Item(11 + $cuisines.Length, 3).ClearContents()
